# Week 12.docx edit:
#  1) In the last row of the "Reports" table, change the "Who Use" cell
#     from "Library Manager" to "Librarian".
#  2) In the same row's "Purpose" cell, change the leading phrase
#     "The library manager would use ..." to "The librarian would use ...".
#  3) Strip the decorative page-number/date drawing shapes out of the
#     footer, leaving the (now empty) footer paragraph in place.

$d = $word.ActiveDocument

# --- 1 & 2: scoped Find/Replace on the specific table cells -------------
# Table has 4 rows; row 4 is "Number of books reserved by each language on
# Feb 2019" - column 4 is "Who Use", column 5 is "Purpose". Using
# $d.Range(start, end) (rather than the cell Range object directly) keeps
# the Find call confined to that cell instead of touching every
# "Library Manager"/"library manager" occurrence in the document.
$t = $d.Tables.Item(1)

$whoUseCell = $t.Cell(4, 4)
$whoUseRange = $d.Range($whoUseCell.Range.Start, $whoUseCell.Range.End)
$whoUseRange.Find.Execute("Library Manager", $false, $false, $false, $false, $false, $true, 0, $false, "Librarian", 1)

$purposeCell = $t.Cell(4, 5)
$purposeRange = $d.Range($purposeCell.Range.Start, $purposeCell.Range.End)
$purposeRange.Find.Execute("The library manager would use", $false, $false, $false, $false, $false, $true, 0, $false, "The librarian would use", 1)

# --- 3: remove the footer's decorative drawing shapes --------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
for ($i = $footer.Shapes.Count; $i -ge 1; $i--) {
    $footer.Shapes.Item($i).Delete()
}
